$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header columns (row 1) ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Title-case the connector words (de/del/el/la/los/y) in municipality / state names ---
$ws.Range("B3").Value = "Pabellón De Arteaga"
$ws.Range("B8").Value = "Amatenango De La Frontera"
$ws.Range("B9").Value = "Amatenango Del Valle"
$ws.Range("A25").Value = "Ciudad De México"
$ws.Range("A37").Value = "Estado De México"
$ws.Range("B40").Value = "Chapa De Mota"
$ws.Range("B43").Value = "Ecatepec De Morelos"
$ws.Range("B45").Value = "Naucalpan De Juárez"
$ws.Range("B48").Value = "San Felipe Del Progreso"
$ws.Range("B51").Value = "Tlalnepantla De Baz"
$ws.Range("B58").Value = "Acapulco De Juárez"
$ws.Range("B60").Value = "Ajuchitlán Del Progreso"
$ws.Range("B63").Value = "Atenango Del Río"
$ws.Range("B64").Value = "Ayutla De Los Libres"
$ws.Range("B65").Value = "Chilapa De Álvarez"
$ws.Range("B66").Value = "Chilpancingo De Los Bravo"
$ws.Range("B68").Value = "Huitzuco De Los Figueroa"
$ws.Range("B69").Value = "Iguala De La Independencia"
$ws.Range("B74").Value = "Técpan De Galeana"
$ws.Range("B75").Value = "Tlapa De Comonfort"
$ws.Range("B80").Value = "Pachuca De Soto"
$ws.Range("B82").Value = "Tulancingo De Bravo"
$ws.Range("B86").Value = "Lagos De Moreno"
$ws.Range("B87").Value = "Tamazula De Gordiano"
$ws.Range("B108").Value = "Chalcatongo De Hidalgo"
$ws.Range("B110").Value = "Eloxochitlán De Flores Magón"
$ws.Range("B111").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B112").Value = "Ixtlán De Juárez"
$ws.Range("B113").Value = "Oaxaca De Juárez"
$ws.Range("B126").Value = "Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca"
$ws.Range("B127").Value = "Zimatlán De Álvarez"
$ws.Range("B130").Value = "Chalchicomula De Sesma"
$ws.Range("B132").Value = "Huehuetlán El Grande"
$ws.Range("B134").Value = "Izúcar De Matamoros"
$ws.Range("B141").Value = "San Salvador El Seco"
$ws.Range("B147").Value = "Tepanco De López"
$ws.Range("B149").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B154").Value = "Ciudad Del Maíz"
$ws.Range("B159").Value = "Soto La Marina"
$ws.Range("B167").Value = "Amatlán De Los Reyes"
$ws.Range("B170").Value = "Las Vigas De Ramírez"
$ws.Range("B172").Value = "Martínez De La Torre"

# --- Remove the trailing metadata/footer rows (185-189 and 476-480) ---
# Delete from the bottom up so row numbers of the remaining rows don't shift
# before we get to them.
$ws.Range("A476:A480").EntireRow.Delete()
$ws.Range("A185:A189").EntireRow.Delete()
